{"js": "// Insert a new paragraph (column break + text) right after the paragraph\n// that ends with \"...s\u1eafp x\u1ebfp theo ng\u00e0y\" (and right before the trailing\n// empty paragraph at the end of the document body).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nconst needle = \"th\u00f4ng tin nguy\u00ean li\u1ec7u v\u00e0 h\u00f3a \u0111\u01a1n nh\u1eadp s\u1eafp x\u1ebfp theo ng\u00e0y\";\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.includes(needle)) {\n    target = paras.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the anchor paragraph for the new note.\");\n}\n\n// The matching run formatting used throughout this document's notes.\nconst rPrXml =\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-US\"/></w:rPr>';\n\n// Office.js's insertBreak() only exposes page/section/line break kinds (no\n// \"column\" kind), so the column break + trailing text run are inserted as\n// raw OOXML (flat-OPC, as required by Range.insertOoxml) to reproduce the\n// exact run/paragraph structure from the authored edit.\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' + rPrXml + '</w:pPr>' +\n  '<w:r>' + rPrXml + '<w:br w:type=\"column\"/></w:r>' +\n  '<w:r>' + rPrXml + '<w:lastRenderedPageBreak/><w:t>B\u1ecf b\u1ea3ng combo g\u1ed9p s\u1ea3n ph\u1ea9m v\u00e0o b\u1ea3ng sanpham</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst afterRange = target.getRange(\"After\");\nafterRange.insertOoxml(flatOpc, \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends with \"...s\u1eafp x\u1ebfp theo ng\u00e0y\" (the anchor for\n# the new note) by searching for its text rather than relying on a\n# hardcoded paragraph index.\n$findRng = $d.Content\n$found = $findRng.Find.Execute(\"th\u00f4ng tin nguy\u00ean li\u1ec7u v\u00e0 h\u00f3a \u0111\u01a1n nh\u1eadp s\u1eafp x\u1ebfp theo ng\u00e0y\")\nif (-not $found) {\n    throw \"Could not find the anchor paragraph for the new note.\"\n}\n$findRng.Expand(4)  # wdParagraph\n$anchorStart = $findRng.Start\n\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $anchorStart) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve the anchor paragraph index.\"\n}\n\n# Make room for the new paragraph right after the anchor (and before the\n# trailing empty paragraph at the end of the document body).\n$anchorRange = $d.Paragraphs.Item($anchorIndex).Range\n$anchorRange.InsertParagraphAfter()\n$target = $d.Paragraphs.Item($anchorIndex + 1).Range\n\n# Word's object model has no \"column break\" WdBreakType constant exposed\n# through a simple Range.InsertBreak call that keeps the break and the\n# following text inside a single paragraph (InsertBreak always mints a\n# fresh paragraph for page/column/section break kinds). Use InsertXML with\n# the exact OOXML run structure instead, so the column break run and the\n# text run land together in the one new paragraph, matching the authored\n# edit.\n$rPr = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-US\"/></w:rPr>'\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr>' + $rPr + '</w:pPr><w:r>' + $rPr + '<w:br w:type=\"column\"/></w:r><w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>B\u1ecf b\u1ea3ng combo g\u1ed9p s\u1ea3n ph\u1ea9m v\u00e0o b\u1ea3ng sanpham</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$target.InsertXML($xml)\n"}
